$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1123.8948
$ws.Range("I32").Value = 959
$ws.Range("J32").Value = 1350.625
$ws.Range("K32").Value = 959
$ws.Range("L32").Value = 1350.625
$ws.Range("M32").Value = -633
$ws.Range("N32").Value = -2002.625
$ws.Range("H107").Value = 930.3333
$ws.Range("I107").Value = 867.8570999999999
$ws.Range("K107").Value = 867.8570999999999
$ws.Range("M107").Value = 1052.1429
$ws.Range("H112").Value = 1926.6
$ws.Range("J112").Value = 2358.375
$ws.Range("L112").Value = 7075.125
$ws.Range("N112").Value = -9291.125
$ws.Range("H116").Value = 2750638.5
$ws.Range("I116").Value = 10991867
$ws.Range("J116").Value = 3562.1904
$ws.Range("K116").Value = 10991867
$ws.Range("L116").Value = 3562.1904
$ws.Range("M116").Value = -10988425
$ws.Range("N116").Value = -10446.1904
$ws.Range("H118").Value = 677.25
$ws.Range("I118").Value = 677.25
$ws.Range("K118").Value = 2031.75
$ws.Range("M118").Value = -374.75
$ws.Range("H135").Value = 397.25
$ws.Range("I135").Value = 357.2
$ws.Range("K135").Value = 3214.8
$ws.Range("M135").Value = -679.7999999999997
$ws.Range("H137").Value = 20919592
$ws.Range("I137").Value = 35716730
$ws.Range("K137").Value = 107150190
$ws.Range("M137").Value = -107147640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1646.25
$ws.Range("I2").Value = 1694.5
$ws.Range("J2").Value = 1501.5
$ws.Range("K2").Value = 1694.5
$ws.Range("L2").Value = 1501.5
$ws.Range("M2").Value = -1581.5
$ws.Range("N2").Value = -1727.5
$ws.Range("H61").Value = 9333
$ws.Range("I61").Value = 4999
$ws.Range("K61").Value = 4999
$ws.Range("M61").Value = -4787
$ws.Range("H116").Value = 1646.25
$ws.Range("I116").Value = 1694.5
$ws.Range("J116").Value = 1501.5
$ws.Range("K116").Value = 1694.5
$ws.Range("L116").Value = 1501.5
$ws.Range("M116").Value = 599.5
$ws.Range("N116").Value = -6089.5
$ws.Range("H123").Value = 41052.668
$ws.Range("J123").Value = 41052.668
$ws.Range("L123").Value = 41052.668
$ws.Range("N123").Value = -50852.668
$ws.Range("H132").Value = 2774.5757
$ws.Range("I132").Value = 2750.3333
$ws.Range("J132").Value = 2817
$ws.Range("K132").Value = 8250.999899999999
$ws.Range("L132").Value = 8451
$ws.Range("M132").Value = -5720.999899999999
$ws.Range("N132").Value = -13511
$ws.Range("H136").Value = 9333
$ws.Range("I136").Value = 4999
$ws.Range("K136").Value = 14997
$ws.Range("M136").Value = -12447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1646.25
$ws.Range("I3").Value = 1694.5
$ws.Range("J3").Value = 1501.5
$ws.Range("K3").Value = 1694.5
$ws.Range("L3").Value = 1501.5
$ws.Range("M3").Value = -1580.5
$ws.Range("N3").Value = -1729.5
$ws.Range("H105").Value = 2952
$ws.Range("I105").Value = 2563.6365
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2563.6365
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -816.6365000000001
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2772.4211
$ws.Range("I16").Value = 1891.1428
$ws.Range("K16").Value = 1891.1428
$ws.Range("M16").Value = -1604.1428
$ws.Range("H22").Value = 708.6667
$ws.Range("I22").Value = 700.6667
$ws.Range("J22").Value = 716.6667
$ws.Range("K22").Value = 700.6667
$ws.Range("L22").Value = 716.6667
$ws.Range("M22").Value = -350.6667
$ws.Range("N22").Value = -1416.6667
$ws.Range("H31").Value = 2795.5945
$ws.Range("I31").Value = 2211.5356
$ws.Range("J31").Value = 4612.6665
$ws.Range("K31").Value = 2211.5356
$ws.Range("L31").Value = 4612.6665
$ws.Range("M31").Value = -1916.5356
$ws.Range("N31").Value = -5202.6665
$ws.Range("H34").Value = 2795.5945
$ws.Range("I34").Value = 2211.5356
$ws.Range("J34").Value = 4612.6665
$ws.Range("K34").Value = 2211.5356
$ws.Range("L34").Value = 4612.6665
$ws.Range("M34").Value = -2009.5356
$ws.Range("N34").Value = -5016.6665
$ws.Range("H113").Value = 2772.4211
$ws.Range("I113").Value = 1891.1428
$ws.Range("K113").Value = 1891.1428
$ws.Range("M113").Value = 278.8571999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 409.1875
$ws.Range("I5").Value = 409.1875
$ws.Range("K5").Value = 1227.5625
$ws.Range("M5").Value = -1115.5625
$ws.Range("H131").Value = 734.9648999999999
$ws.Range("J131").Value = 954.3714
$ws.Range("L131").Value = 2863.1142
$ws.Range("N131").Value = -12943.1142
$ws.Range("H135").Value = 409.1875
$ws.Range("I135").Value = 409.1875
$ws.Range("K135").Value = 3682.6875
$ws.Range("M135").Value = -1147.6875
$ws.Range("H136").Value = 3822.5
$ws.Range("I136").Value = 1895
$ws.Range("K136").Value = 5685
$ws.Range("M136").Value = -585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2935.8
$ws.Range("I126").Value = 2178
$ws.Range("K126").Value = 6534
$ws.Range("M126").Value = -4064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 35894
$ws.Range("J87").Value = 35894
$ws.Range("L87").Value = 35894
$ws.Range("N87").Value = -38140
$ws.Range("H90").Value = 35894
$ws.Range("J90").Value = 35894
$ws.Range("L90").Value = 107682
$ws.Range("N90").Value = -118914
$ws.Range("H132").Value = 2887.8
$ws.Range("I132").Value = 2357.3125
$ws.Range("K132").Value = 7071.9375
$ws.Range("M132").Value = -4541.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 745.125
$ws.Range("I81").Value = 737.2857
$ws.Range("J81").Value = 800
$ws.Range("K81").Value = 1474.5714
$ws.Range("L81").Value = 1600
$ws.Range("M81").Value = -413.5714
$ws.Range("N81").Value = -3722
$ws.Range("H84").Value = 745.125
$ws.Range("I84").Value = 737.2857
$ws.Range("J84").Value = 800
$ws.Range("K84").Value = 7372.857
$ws.Range("L84").Value = 8000
$ws.Range("M84").Value = -2068.857
$ws.Range("N84").Value = -18608
$ws.Range("H107").Value = 704
$ws.Range("I107").Value = 440.63635
$ws.Range("J107").Value = 1066.125
$ws.Range("K107").Value = 1321.90905
$ws.Range("L107").Value = 3198.375
$ws.Range("M107").Value = 598.09095
$ws.Range("N107").Value = -7038.375
